# "Updates, Small bug fixes"
#
# Nudges the positions of the species-label textboxes (tx9-tx18) that sit
# inside the chart group on slide 1. Only each shape's offset (Left/Top)
# moves a little; the sizes (Width/Height) are unchanged.
#
# Shape.Left / Shape.Top are exposed by PowerPoint as points (1 pt =
# 12700 EMU) but are stored internally as single-precision (float32)
# values, so the literals below are the exact float32 values whose
# point->EMU conversion reproduces the target EMU offsets bit-for-bit
# (rather than the naive division, which can land one EMU off after the
# float32 round-trip).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The labels live inside the single group shape on the slide.
$grp = $s.Shapes.Item(2)

# Shape name -> new (Left, Top) in points, pre-compensated for the
# float32 storage so the saved EMU offsets match the target exactly.
$moves = @{
    "tx9"  = @(449.2007141113281, 338.2904052734375)
    "tx10" = @(472.9118347167969, 366.05780029296875)
    "tx11" = @(401.1317443847656, 432.2406311035156)
    "tx12" = @(415.8319091796875, 456.60711669921875)
    "tx13" = @(301.0242614746094, 370.7319030761719)
    "tx14" = @(324.74371337890625, 395.0983581542969)
    "tx15" = @(320.5530090332031, 283.4993896484375)
    "tx16" = @(325.2877197265625, 307.8658447265625)
    "tx17" = @(345.14678955078125, 240.0255126953125)
    "tx18" = @(349.4104919433594, 264.3836364746094)
}

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $item = $grp.GroupItems.Item($i)
    if ($moves.ContainsKey($item.Name)) {
        $xy = $moves[$item.Name]
        $item.Left = $xy[0]
        $item.Top  = $xy[1]
    }
}
